$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0,0,1,0,1,0,0,1,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("C2").Select()
